# Update the cached "datetimeFigureOut" date field text from 2018/1/12 to
# 2018/1/28 on every date placeholder across the slide master and all of
# its slide layouts (the ROS / ROS-lite figure masters/layouts).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        if ($shp.Type -eq 14) {
            # msoPlaceholder
            if ($shp.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        }

        $tf = $shp.TextFrame
        if (-not $tf.HasText) { continue }
        $txt = $tf.TextRange.Text

        if ($isDatePlaceholder -or ($txt -eq "2018/1/12")) {
            if ($txt -eq "2018/1/12") {
                $tf.TextRange.Text = "2018/1/28"
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder($p.SlideMaster)

# Every custom layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder($layouts.Item($L))
}
